# Auto-generated edit script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.806.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.20%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.488.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.33%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'592.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.11%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'171.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.35%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.06%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.131"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.46%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -1.24%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.37%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'4.088.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.36%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -0.44%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'28.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.95%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'66.806.79"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "'  -1.48%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.483.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.85%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'6.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.26%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'14.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.12%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'392.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.21%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'7.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.65%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'72.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.85%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.62%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.0000120"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.54%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'10.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.76%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -0.52%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -0.02%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'6.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.11%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.66%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -1.46%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'23.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.07%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'7.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.82%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -0.86%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'162.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.07%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.877"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.65%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.74%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'6.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.28%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'4.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.46%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0739"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.24%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'27.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.86%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'26.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.81%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.795.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.14%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'42.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.29%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +1.59%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -3.60%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'337.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.01%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'34.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.71%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -3.04%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -1.58%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'6.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.52%  "
$ws.Range("E51").Style = "Normal"
